# Auto-generated edit script applying cryptos.xlsx diff
# (cryptos list refreshed by the GitHub Actions scraper)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the value to be written as plain text, exactly as the
    # source data supplies it (avoids Excel auto-converting numeric-
    # looking strings like '2.30' or '7.50' into numbers and silently
    # dropping the trailing zero), then restores the default style so
    # no stray cell formatting is introduced.
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.Style = 'Normal'
}

Set-TextCell $ws.Range('D2') '26.663.74'
Set-TextCell $ws.Range('E2') '  -1.47%  '
Set-TextCell $ws.Range('D3') '1.595.93'
Set-TextCell $ws.Range('E3') '  -1.62%  '
Set-TextCell $ws.Range('E4') '  +0.03%  '
Set-TextCell $ws.Range('D5') '211.08'
Set-TextCell $ws.Range('E5') '  -1.31%  '
Set-TextCell $ws.Range('D6') '0.511'
Set-TextCell $ws.Range('E6') '  -0.12%  '
Set-TextCell $ws.Range('E7') '  +0.04%  '
Set-TextCell $ws.Range('E8') '  -1.37%  '
Set-TextCell $ws.Range('E9') '  -1.52%  '
Set-TextCell $ws.Range('D10') '19.66'
Set-TextCell $ws.Range('E10') '  -1.32%  '
Set-TextCell $ws.Range('E11') '  -0.50%  '
Set-TextCell $ws.Range('D12') '1.818.19'
Set-TextCell $ws.Range('E12') '  -1.67%  '
Set-TextCell $ws.Range('D13') '1.591.95'
Set-TextCell $ws.Range('E13') '  -1.97%  '
Set-TextCell $ws.Range('E14') '  -2.29%  '
Set-TextCell $ws.Range('D15') '0.523'
Set-TextCell $ws.Range('E15') '  -2.99%  '
Set-TextCell $ws.Range('D16') '64.83'
Set-TextCell $ws.Range('E16') '  +0.57%  '
Set-TextCell $ws.Range('D17') '26.638.37'
Set-TextCell $ws.Range('E17') '  -1.52%  '
Set-TextCell $ws.Range('D18') '0.0₃0730'
Set-TextCell $ws.Range('E18') '  -0.88%  '
Set-TextCell $ws.Range('D19') '209.14'
Set-TextCell $ws.Range('E19') '  -2.55%  '
Set-TextCell $ws.Range('E20') '  +0.05%  '
Set-TextCell $ws.Range('E21') '  -2.12%  '
Set-TextCell $ws.Range('E22') '  -2.08%  '
Set-TextCell $ws.Range('D23') '2.30'
Set-TextCell $ws.Range('E23') '  -0.82%  '
Set-TextCell $ws.Range('D24') '8.90'
Set-TextCell $ws.Range('E24') '  -1.21%  '
Set-TextCell $ws.Range('D25') '146.62'
Set-TextCell $ws.Range('E25') '  -0.49%  '
Set-TextCell $ws.Range('E26') '  +0.06%  '
Set-TextCell $ws.Range('D27') '7.12'
Set-TextCell $ws.Range('E27') '  -4.18%  '
Set-TextCell $ws.Range('E28') '  -0.04%  '
Set-TextCell $ws.Range('D29') '15.31'
Set-TextCell $ws.Range('E29') '  -1.28%  '
Set-TextCell $ws.Range('E30') '  -1.60%  '
Set-TextCell $ws.Range('E31') '  -1.20%  '
Set-TextCell $ws.Range('E32') '  -2.71%  '
Set-TextCell $ws.Range('D33') '0.686'
Set-TextCell $ws.Range('E33') '  -4.36%  '
Set-TextCell $ws.Range('E34') '  -2.88%  '
Set-TextCell $ws.Range('D35') '1.293.19'
Set-TextCell $ws.Range('E35') '  -3.27%  '
Set-TextCell $ws.Range('E36') '  -0.58%  '
Set-TextCell $ws.Range('E37') '  -4.89%  '
Set-TextCell $ws.Range('E38') '  -2.64%  '
Set-TextCell $ws.Range('D39') '0.836'
Set-TextCell $ws.Range('E39') '  -0.27%  '
Set-TextCell $ws.Range('E40') '  +0.05%  '
Set-TextCell $ws.Range('E41') '  -0.33%  '
Set-TextCell $ws.Range('B42') 'FraxShare'
Set-TextCell $ws.Range('C42') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws.Range('D42') '5.35'
Set-TextCell $ws.Range('E42') '  +0.18%  '
Set-TextCell $ws.Range('B43') 'MXToken'
Set-TextCell $ws.Range('C43') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws.Range('D43') '2.19'
Set-TextCell $ws.Range('E43') '  -1.66%  '
Set-TextCell $ws.Range('D44') '63.55'
Set-TextCell $ws.Range('E44') '  -0.51%  '
Set-TextCell $ws.Range('D45') '1.730.56'
Set-TextCell $ws.Range('E45') '  -1.70%  '
Set-TextCell $ws.Range('D46') '0.897'
Set-TextCell $ws.Range('E46') '  +4.57%  '
Set-TextCell $ws.Range('E47') '  -0.23%  '
Set-TextCell $ws.Range('D48') '1.64'
Set-TextCell $ws.Range('E48') '  -0.17%  '
Set-TextCell $ws.Range('B49') 'Algorand'
Set-TextCell $ws.Range('C49') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws.Range('D49') '0.0983'
Set-TextCell $ws.Range('E49') '  -1.29%  '
Set-TextCell $ws.Range('B50') 'Cronos'
Set-TextCell $ws.Range('C50') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range('D50') '0.0503'
Set-TextCell $ws.Range('E50') '  -1.70%  '
Set-TextCell $ws.Range('B51') 'EnergySwap'
Set-TextCell $ws.Range('D51') '7.50'
Set-TextCell $ws.Range('E51') '  -0.83%  '
